$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the red font style first so it becomes cellXfs index 4 (matches target order)
$ws.Range("A36").Font.Color = 255

# Row 9 (Robustness) - add D9 = 1 with yellow highlight
$ws.Range("D9").Value = 1
$ws.Range("D9").Interior.Color = 65535

# Row 26 (Operator Overloading) - add D26 = 1 with yellow highlight
$ws.Range("D26").Value = 1
$ws.Range("D26").Interior.Color = 65535

# Row 27 (Inheritance) - add D27 = 1
$ws.Range("D27").Value = 1

# Row 28 (Virtual Function / Overriding) - add D28 = 2
$ws.Range("D28").Value = 2

# Row 31 (Bitwise Operators) - add D31 = 1 with yellow highlight
$ws.Range("D31").Value = 1
$ws.Range("D31").Interior.Color = 65535

# Row 32 (Concurrent Programming) - add D32 = 1 with yellow highlight
$ws.Range("D32").Value = 1
$ws.Range("D32").Interior.Color = 65535

# D33 total (SUM(D3:D32)) recalculates automatically from the new values above

# Update the view: scroll position and selection
$ws.Range("D20").Select() | Out-Null
